$d = $word.ActiveDocument

# 1) The third empty paragraph (right before the "NPAPI:" paragraph) gets a
#    StackOverflow link typed into it.
$target = $d.Paragraphs.Item(4)
$target.Range.InsertAfter("https://stackoverflow.com/questions/42288596/websql-has-increasing-browser-support-whats-its-future")

# 2) Remove the stray <w:lastRenderedPageBreak/> that sits just before the
#    "Unterteilung in local storage..." text (between the localstorage.html
#    hyperlink and that text there are three manual line breaks; the middle
#    one carries the stale lastRenderedPageBreak marker).
$find = $d.Content
$find.Find.ClearFormatting()
$ok = $find.Find.Execute("Unterteilung in local storage")
if ($ok) {
    $brPos = $find.Start - 2
    $d.Range($brPos, $brPos + 1).Delete()
}

# 3) Move the "_GoBack" bookmark from the end of the document onto the newly
#    typed StackOverflow paragraph (re-adding a bookmark with the same name
#    relocates it, since bookmark names are unique in a document).
$newPara = $d.Paragraphs.Item(4).Range
$bmRange = $d.Range($newPara.Start, $newPara.End - 1)
$d.Bookmarks.Add("_GoBack", $bmRange)
